# Refactor synthetic array /3: recolor status indicators from black/red/green/orange
# squares to blue/red/green/orange books, and rename "noir" (black) to "bleu" (blue).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange

# Map of old value -> new value (emoji square indicators -> emoji book indicators,
# and the matching color label).
$replacements = @{
    "⬛"   = "📘"
    "🟥"   = "📕"
    "🟩"   = "📗"
    "🟧"   = "📙"
    "noir" = "bleu"
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $usedRange.Replace($old, $new, 1, 1, $false, $false, $false, $false) | Out-Null
}
